$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (stored OOXML "width" = ColumnWidth + 5/6 padding,
# so subtract the padding to land on the target stored widths).
$ws.Range("A1").ColumnWidth = 1.3307291666666665   # -> 2.1640625
$ws.Range("B1").ColumnWidth = 2.4986979166666665   # -> 3.33203125
$ws.Range("C1").ColumnWidth = 11.998697916666666   # -> 12.83203125
$ws.Range("D1").ColumnWidth = 18.830729166666668   # -> 19.6640625
$ws.Range("E1").ColumnWidth = 5.330729166666667    # -> 6.1640625
$ws.Range("F1").ColumnWidth = 2.4986979166666665   # -> 3.33203125
$ws.Range("G1").ColumnWidth = 1.3307291666666665   # -> 2.1640625
$ws.Range("H1").ColumnWidth = 0.9986979166666666   # -> 1.83203125
$ws.Range("I1").ColumnWidth = 11.998697916666666   # -> 12.83203125
$ws.Range("J1").ColumnWidth = 5.330729166666667    # -> 6.1640625
$ws.Range("K1").ColumnWidth = 0.9986979166666666   # -> 1.83203125
$ws.Range("L1").ColumnWidth = 1.3307291666666665   # -> 2.1640625
$ws.Range("M1").ColumnWidth = 9.666666666666666    # -> 10.5
$ws.Range("O1").ColumnWidth = 34.666666666666664   # -> 35.5
$ws.Range("T1").ColumnWidth = 0.8307291666666666   # -> 1.6640625

# Move selection to I9 (also resets the scrolled topLeftCell back to the
# sheet's natural top-left).
$ws.Range("I9").Select()
